$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3081.818
$ws.Range("I51").Value = 1700
$ws.Range("J51").Value = 3600
$ws.Range("K51").Value = 1700
$ws.Range("L51").Value = 3600
$ws.Range("M51").Value = -1216
$ws.Range("N51").Value = -4568
$ws.Range("H132").Value = 4170093.2
$ws.Range("I132").Value = 3318.6155
$ws.Range("J132").Value = 31254128
$ws.Range("K132").Value = 9955.8465
$ws.Range("L132").Value = 93762384
$ws.Range("M132").Value = -7425.8465
$ws.Range("N132").Value = -93767444
$ws.Range("H138").Value = 7248202.5
$ws.Range("I138").Value = 12346873
$ws.Range("J138").Value = 2723.6843
$ws.Range("K138").Value = 37040619
$ws.Range("L138").Value = 8171.0529
$ws.Range("M138").Value = -37035479
$ws.Range("N138").Value = -18451.0529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10271.974
$ws.Range("I32").Value = 10840.121
$ws.Range("J32").Value = 8537.632
$ws.Range("K32").Value = 10840.121
$ws.Range("L32").Value = 8537.632
$ws.Range("M32").Value = -10553.121
$ws.Range("N32").Value = -9111.632
$ws.Range("H45").Value = 1431.4482
$ws.Range("I45").Value = 1305.091
$ws.Range("J45").Value = 1828.5714
$ws.Range("K45").Value = 1305.091
$ws.Range("L45").Value = 1828.5714
$ws.Range("M45").Value = -928.0909999999999
$ws.Range("N45").Value = -2582.5714
$ws.Range("H122").Value = 5131.727
$ws.Range("I122").Value = 6149.4585
$ws.Range("J122").Value = 2417.7778
$ws.Range("K122").Value = 18448.3755
$ws.Range("L122").Value = 7253.3334
$ws.Range("M122").Value = -15998.3755

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1950.174
$ws.Range("I105").Value = 999.13336
$ws.Range("J105").Value = 3733.375
$ws.Range("K105").Value = 999.13336
$ws.Range("L105").Value = 3733.375
$ws.Range("M105").Value = 747.86664
$ws.Range("H134").Value = 2365.0833
$ws.Range("I134").Value = 1413.9722
$ws.Range("J134").Value = 3791.75
$ws.Range("K134").Value = 4241.9166
$ws.Range("L134").Value = 11375.25
$ws.Range("M134").Value = -1706.9166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 41679.8
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 41679.8
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 41679.8
$ws.Range("N20").Value = -42151.8
$ws.Range("H30").Value = 41679.8
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 41679.8
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 41679.8
$ws.Range("N30").Value = -41861.8
$ws.Range("H58").Value = 1674.3334
$ws.Range("I58").Value = 774.3077
$ws.Range("J58").Value = 2510.0715
$ws.Range("K58").Value = 774.3077
$ws.Range("L58").Value = 2510.0715
$ws.Range("M58").Value = -571.3077
$ws.Range("N58").Value = -2916.0715
$ws.Range("H86").Value = 3619.875
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 3391.8
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 3391.8
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -5637.8
$ws.Range("H89").Value = 3619.875
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 3391.8
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 16959
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -28191
$ws.Range("H99").Value = 1464.7858
$ws.Range("I99").Value = 1500.5385
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1500.5385
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -2.538500000000113
$ws.Range("N99").Value = -3996
$ws.Range("H122").Value = 1473.95
$ws.Range("I122").Value = 1665.625
$ws.Range("J122").Value = 707.25
$ws.Range("K122").Value = 4996.875
$ws.Range("L122").Value = 2121.75
$ws.Range("M122").Value = -2546.875
$ws.Range("N122").Value = -7021.75
$ws.Range("H126").Value = 1464.7858
$ws.Range("I126").Value = 1500.5385
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 4501.6155
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -2031.6155
$ws.Range("N126").Value = -7940
$ws.Range("H128").Value = 41679.8
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41679.8
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41679.8
$ws.Range("N128").Value = -51639.8
$ws.Range("H132").Value = 8930257
$ws.Range("I132").Value = 11906281
$ws.Range("J132").Value = 2187.1428
$ws.Range("K132").Value = 35718843
$ws.Range("L132").Value = 6561.428400000001
$ws.Range("M132").Value = -35716313
$ws.Range("N132").Value = -11621.4284
$ws.Range("H134").Value = 373238.2
$ws.Range("I134").Value = 1276.7609
$ws.Range("J134").Value = 1323806.2
$ws.Range("K134").Value = 3830.2827
$ws.Range("L134").Value = 3971418.6
$ws.Range("M134").Value = -1295.2827
$ws.Range("N134").Value = -3976488.6
$ws.Range("H136").Value = 1674.3334
$ws.Range("I136").Value = 774.3077
$ws.Range("J136").Value = 2510.0715
$ws.Range("K136").Value = 2322.9231
$ws.Range("L136").Value = 7530.2145
$ws.Range("M136").Value = 227.0769
$ws.Range("N136").Value = -12630.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 708.4857
$ws.Range("I5").Value = 259.18182
$ws.Range("J5").Value = 1468.8462
$ws.Range("K5").Value = 777.54546
$ws.Range("L5").Value = 4406.5386
$ws.Range("M5").Value = -665.54546
$ws.Range("N5").Value = -4630.5386
$ws.Range("H122").Value = 848.381
$ws.Range("I122").Value = 1493.375
$ws.Range("J122").Value = 451.46155
$ws.Range("K122").Value = 13440.375
$ws.Range("L122").Value = 4063.15395
$ws.Range("M122").Value = -10990.375
$ws.Range("N122").Value = -8963.15395
$ws.Range("H135").Value = 708.4857
$ws.Range("I135").Value = 259.18182
$ws.Range("J135").Value = 1468.8462
$ws.Range("K135").Value = 2332.63638
$ws.Range("L135").Value = 13219.6158
$ws.Range("M135").Value = 202.3636200000001
$ws.Range("N135").Value = -18289.6158
$ws.Range("H136").Value = 2428.4546
$ws.Range("I136").Value = 1682.8572
$ws.Range("J136").Value = 3733.25
$ws.Range("K136").Value = 5048.571599999999
$ws.Range("L136").Value = 11199.75
$ws.Range("M136").Value = 51.42840000000069
$ws.Range("N136").Value = -21399.75
$ws.Range("H139").Value = 3218.4211
$ws.Range("I139").Value = 1695.4546
$ws.Range("J139").Value = 5312.5
$ws.Range("K139").Value = 5086.3638
$ws.Range("L139").Value = 15937.5
$ws.Range("M139").Value = 53.63619999999992
$ws.Range("H141").Value = 7186.364
$ws.Range("I141").Value = 5505.5557
$ws.Range("J141").Value = 14750
$ws.Range("K141").Value = 16516.6671
$ws.Range("L141").Value = 44250
$ws.Range("M141").Value = -11336.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13212.956
$ws.Range("I70").Value = 26200.889
$ws.Range("J70").Value = 4863.5713
$ws.Range("K70").Value = 26200.889
$ws.Range("L70").Value = 4863.5713
$ws.Range("M70").Value = -25930.889
$ws.Range("N70").Value = -5403.5713
$ws.Range("H73").Value = 13212.956
$ws.Range("I73").Value = 26200.889
$ws.Range("J73").Value = 4863.5713
$ws.Range("K73").Value = 26200.889
$ws.Range("L73").Value = 4863.5713
$ws.Range("M73").Value = -25264.889
$ws.Range("N73").Value = -6735.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1073.6154
$ws.Range("I46").Value = 684
$ws.Range("J46").Value = 1317.125
$ws.Range("K46").Value = 684
$ws.Range("L46").Value = 1317.125
$ws.Range("M46").Value = -496
$ws.Range("N46").Value = -1693.125
$ws.Range("H68").Value = 1652
$ws.Range("I68").Value = 1697
$ws.Range("J68").Value = 1580
$ws.Range("K68").Value = 1697
$ws.Range("L68").Value = 1580
$ws.Range("M68").Value = -948
$ws.Range("N68").Value = -3078
$ws.Range("H71").Value = 1652
$ws.Range("I71").Value = 1697
$ws.Range("J71").Value = 1580
$ws.Range("K71").Value = 8485
$ws.Range("L71").Value = 7900
$ws.Range("M71").Value = -4741
$ws.Range("N71").Value = -15388

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1396
$ws.Range("I81").Value = 1226.6666
$ws.Range("J81").Value = 1468.5714
$ws.Range("K81").Value = 2453.3332
$ws.Range("L81").Value = 2937.1428
$ws.Range("M81").Value = -1392.3332
$ws.Range("N81").Value = -5059.1428
$ws.Range("H84").Value = 1396
$ws.Range("I84").Value = 1226.6666
$ws.Range("J84").Value = 1468.5714
$ws.Range("K84").Value = 12266.666
$ws.Range("L84").Value = 14685.714
$ws.Range("M84").Value = -6962.666000000001
$ws.Range("N84").Value = -25293.714
$ws.Range("H107").Value = 620.13794
$ws.Range("I107").Value = 657.36
$ws.Range("J107").Value = 387.5
$ws.Range("K107").Value = 1972.08
$ws.Range("L107").Value = 1162.5
$ws.Range("M107").Value = -52.07999999999993
$ws.Range("N107").Value = -5002.5
